# Update the "Metadata" sheet (StructureDefinition property/value table).
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

# Version bump 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> now populated
$meta.Range("B9").Value = "Alvearie Team"

# Remove the duplicated "Contact" / "No display for ContactDetail" row (row 11),
# which shifts every row below it up by one.
$meta.Rows.Item(11).Delete()

# The row that used to hold the second "Contact" entry (row 10) becomes the new
# "Jurisdiction" property.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Update the "Elements" sheet: the root Extension row's Short/Definition columns
# now carry the profile-specific title/description instead of the generic
# "Extension" / "An Extension" placeholders.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Mental Health Day Night Coverage Indicator"
$elements.Range("L2").Value = "Indicates whether the member has mental health day-night benefit coverage: Y or N. This finer granularity of MHSA benefit coverage may be used in HEDIS reporting."

# Widen column K ("Short") on the Elements sheet to fit the new (longer) text
# (target stored width ~41.63 chars; 40.8 is the closest settable value given
# this host's internal pixel-width quantization).
$elements.Columns.Item(11).ColumnWidth = 40.8
